# Update the "Estado de Cuenta" worker table:
#  - Reorder the worker rows (16-20): LEONOR/ZULEY swap, and CRISTINA/LICETH swap
#  - Update "Periodo Mora" column (E) from 2506 to 2507 for all worker rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the worker table (rows 16-20), columns C (Doc #), D (Name)
$ws.Range("C16").Value = "45505666"
$ws.Range("D16").Value = "LEONOR PATRICIA MIER GOMEZ"

$ws.Range("C17").Value = "1047475467"
$ws.Range("D17").Value = "ZULEY VASQUEZ CANABAL"

$ws.Range("C18").Value = "1047413174"
$ws.Range("D18").Value = "MATIS JOHANA ORTEGA PALOMINO"

$ws.Range("C19").Value = "1047448808"
$ws.Range("D19").Value = "LICETH PAOLA AYALA HERNANDEZ"

$ws.Range("C20").Value = "1047468999"
$ws.Range("D20").Value = "CRISTINA PACHECO ALVAREZ"

# Update Periodo Mora (column E) for all worker rows to 2507
$ws.Range("E16:E20").Value = "2507"
